$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "manage address" rows appended below existing data
$ws.Range("A7").Value = "name"
$ws.Range("B7").Value = "Naruto Uzumaki"

$ws.Range("A8").Value = "number"
$ws.Range("B8").Value = 9999999106

$ws.Range("A9").Value = "pincode"
$ws.Range("B9").Value = 500085

$ws.Range("A10").Value = "locality"
$ws.Range("B10").Value = "Konoha"

$ws.Range("A11").Value = "address"
$ws.Range("B11").Value = "D-no 11/a, Uzumaki Street, Hidden Leaf Village, Konoha"

# Column B width widened to fit new content (stored sheet width ends up
# offset from the COM ColumnWidth by Excel's default column-padding, so
# compensate to land on a stored width of 14.5)
$ws.Columns.Item(2).ColumnWidth = 13.6666666666667

# Update selection to match final state
$ws.Range("A10").Select()
